$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.033305856709975
$ws.Range("D2").Value = 1.037653753515076
$ws.Range("E2").Value = 0.9926147277508489
$ws.Range("F2").Value = 1.031939001865432
$ws.Range("I2").Value = 1.037995726329598
$ws.Range("J2").Value = 1.038431309000492
$ws.Range("K2").Value = 1.040444022575654
$ws.Range("L2").Value = 0.9955398523336033
$ws.Range("M2").Value = 1.034745672920536
$ws.Range("N2").Value = 1.016716152168733
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.034416717210054
$ws.Range("D3").Value = 1.038521502404954
$ws.Range("E3").Value = 0.9936372048519304
$ws.Range("F3").Value = 1.033681716798242
$ws.Range("I3").Value = 1.038349306449945
$ws.Range("J3").Value = 1.039184009651757
$ws.Range("K3").Value = 1.04112147028148
$ws.Range("L3").Value = 0.9963617723202692
$ws.Range("M3").Value = 1.036294556829217
$ws.Range("N3").Value = 1.016971395930734
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.035134342349869
$ws.Range("D4").Value = 1.039081898456121
$ws.Range("E4").Value = 0.9942998659930995
$ws.Range("F4").Value = 1.034808048075874
$ws.Range("I4").Value = 1.038576106959262
$ws.Range("J4").Value = 1.039669350303737
$ws.Range("K4").Value = 1.04155809425172
$ws.Range("L4").Value = 0.9968940712668345
$ws.Range("M4").Value = 1.037295030839944
$ws.Range("N4").Value = 1.017135843460491
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.035435753766698
$ws.Range("D5").Value = 1.039317228672967
$ws.Range("E5").Value = 0.9945786998346017
$ws.Range("F5").Value = 1.035281251253029
$ws.Range("I5").Value = 1.038670978904319
$ws.Range("J5").Value = 1.039872981343277
$ws.Range("K5").Value = 1.041741238849853
$ws.Range("L5").Value = 0.997117960005301
$ws.Range("M5").Value = 1.037715218207485
$ws.Range("N5").Value = 1.017204807485787
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.035486345852499
$ws.Range("D6").Value = 1.039356726445964
$ws.Range("E6").Value = 0.9946255319796338
$ws.Range("F6").Value = 1.035360686464047
$ws.Range("I6").Value = 1.038686880505135
$ws.Range("J6").Value = 1.039907148144727
$ws.Range("K6").Value = 1.04177196552596
$ws.Range("L6").Value = 0.9971555583673453
$ws.Range("M6").Value = 1.037785745616671
$ws.Range("N6").Value = 1.017216376919929
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.03513837091298
$ws.Range("D7").Value = 1.039085043970303
$ws.Range("E7").Value = 0.9943035907982488
$ws.Range("F7").Value = 1.034814372231892
$ws.Range("I7").Value = 1.038577376507634
$ws.Range("J7").Value = 1.039672072824567
$ws.Range("K7").Value = 1.041560543055642
$ws.Range("L7").Value = 0.9968970624462087
$ws.Range("M7").Value = 1.037300647011725
$ws.Range("N7").Value = 1.017136765627056
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.033681522695562
$ws.Range("D8").Value = 1.037947241347132
$ws.Range("E8").Value = 0.9929600610674301
$ws.Range("F8").Value = 1.032528240233354
$ws.Range("I8").Value = 1.038115633102962
$ws.Range("J8").Value = 1.038686042659904
$ws.Range("K8").Value = 1.040673328643568
$ws.Range("L8").Value = 0.995817528259106
$ws.Range("M8").Value = 1.035269494709713
$ws.Range("N8").Value = 1.016802560952063
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.031105214372021
$ws.Range("D9").Value = 1.035933807374396
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.0284892128835
$ws.Range("I9").Value = 1.037286686531684
$ws.Range("J9").Value = 1.036935345386832
$ws.Range("K9").Value = 1.039096601096069
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.031676488067223
$ws.Range("N9").Value = 1.016208158506152
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.029381307295348
$ws.Range("D10").Value = 1.034585687035092
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.025788793877226
$ws.Range("I10").Value = 1.036723687560284
$ws.Range("J10").Value = 1.035759188245554
$ws.Range("K10").Value = 1.038036346405276
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.029271277050404
$ws.Range("N10").Value = 1.015808146868951
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.028633277833727
$ws.Range("D11").Value = 1.034000523227245
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.024617508194593
$ws.Range("I11").Value = 1.036477424392942
$ws.Range("J11").Value = 1.035247722047239
$ws.Range("K11").Value = 1.037575055037318
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.028227327885361
$ws.Range("N11").Value = 1.015634037807649
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.028355186798863
$ws.Range("D12").Value = 1.033782951497556
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.024182129973621
$ws.Range("I12").Value = 1.036385576821115
$ws.Range("J12").Value = 1.035057409610436
$ws.Range("K12").Value = 1.037403378536868
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.02783917575818
$ws.Range("N12").Value = 1.015569229478372
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.028414849180535
$ws.Range("D13").Value = 1.033829631182558
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.024275534329493
$ws.Range("I13").Value = 1.03640529541059
$ws.Range("J13").Value = 1.035098247310259
$ws.Range("K13").Value = 1.03744021881223
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.0279224532095
$ws.Range("N13").Value = 1.015583137277585
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.028610295654337
$ws.Range("D14").Value = 1.033982543098813
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.024581526100999
$ws.Range("I14").Value = 1.036469839893138
$ws.Range("J14").Value = 1.035231997535388
$ws.Range("K14").Value = 1.037560870999803
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.028195250983763
$ws.Range("N14").Value = 1.01562868352102
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.028730684801585
$ws.Range("D15").Value = 1.034076728495985
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.024770016212842
$ws.Range("I15").Value = 1.036509558206249
$ws.Range("J15").Value = 1.035314361489692
$ws.Range("K15").Value = 1.03763516467376
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.02836327963018
$ws.Range("N15").Value = 1.015656727949433
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.029430918183405
$ws.Range("D16").Value = 1.034624492309434
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.025866485365589
$ws.Range("I16").Value = 1.036739978817916
$ws.Range("J16").Value = 1.035793086316502
$ws.Range("K16").Value = 1.038066914336527
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.02934050741288
$ws.Range("N16").Value = 1.015819682833519
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.029869733735916
$ws.Range("D17").Value = 1.034967708343003
$ws.Range("E17").Value = 0.9894763578477731
$ws.Range("F17").Value = 1.026553731139904
$ws.Range("I17").Value = 1.036883850150624
$ws.Range("J17").Value = 1.036092791187647
$ws.Range("K17").Value = 1.038337150168648
$ws.Range("L17").Value = 0.9930127773692701
$ws.Range("M17").Value = 1.029952825845284
$ws.Range("N17").Value = 1.015921658063723
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.030125536587242
$ws.Range("D18").Value = 1.035167763738083
$ws.Range("E18").Value = 0.9897087662937551
$ws.Range("F18").Value = 1.026954399135216
$ws.Range("I18").Value = 1.036967528613244
$ws.Range("J18").Value = 1.036267393620047
$ws.Range("K18").Value = 1.038494562644961
$ws.Range("L18").Value = 0.9932001317071766
$ws.Range("M18").Value = 1.030309742574142
$ws.Range("N18").Value = 1.015981051591528
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.030212733292611
$ws.Range("D19").Value = 1.035235954349571
$ws.Range("E19").Value = 0.9897880325774039
$ws.Range("F19").Value = 1.027090984563341
$ws.Range("I19").Value = 1.036996020277476
$ws.Range("J19").Value = 1.036326892955799
$ws.Range("K19").Value = 1.03854820043653
$ws.Range("L19").Value = 0.993264023964098
$ws.Range("M19").Value = 1.030431401868368
$ws.Range("N19").Value = 1.016001288541592
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.029822668613871
$ws.Range("D20").Value = 1.034930898671711
$ws.Range("E20").Value = 0.9894336180360677
$ws.Range("F20").Value = 1.026480016001666
$ws.Range("I20").Value = 1.03686843887829
$ws.Range("J20").Value = 1.036060657473058
$ws.Range("K20").Value = 1.038308178319562
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.029887154596309
$ws.Range("N20").Value = 1.015910726082869
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.028552748187388
$ws.Range("D21").Value = 1.033937520330531
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.024491427821979
$ws.Range("I21").Value = 1.036450843502933
$ws.Range("J21").Value = 1.035192620594355
$ws.Range("K21").Value = 1.037525351134067
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.028114929496292
$ws.Range("N21").Value = 1.015615275063302
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.027752910497369
$ws.Range("D22").Value = 1.033311693824298
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.023239320586776
$ws.Range("I22").Value = 1.036186117292852
$ws.Range("J22").Value = 1.034644933984992
$ws.Range("K22").Value = 1.037031232389386
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.026998440339686
$ws.Range("N22").Value = 1.015428723170841
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.028177052621146
$ws.Range("D23").Value = 1.033643575695893
$ws.Range("E23").Value = 0.9879432794643023
$ws.Range("F23").Value = 1.023903261295188
$ws.Range("I23").Value = 1.036326659698365
$ws.Range("J23").Value = 1.0349354558538
$ws.Range("K23").Value = 1.037293357392137
$ws.Range("L23").Value = 0.991776070289318
$ws.Range("M23").Value = 1.027590526510456
$ws.Range("N23").Value = 1.015527693105143
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.029843935782455
$ws.Range("D24").Value = 1.034947531800023
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.02651332528816
$ws.Range("I24").Value = 1.036875403308039
$ws.Range("J24").Value = 1.03607517796525
$ws.Range("K24").Value = 1.038321270104479
$ws.Range("L24").Value = 0.9929938892766442
$ws.Range("M24").Value = 1.029916829342122
$ws.Range("N24").Value = 1.015915666042904
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.031772358571091
$ws.Range("D25").Value = 1.036455346491991
$ws.Range("E25").Value = 0.9912096547607049
$ws.Range("F25").Value = 1.029534715221714
$ws.Range("I25").Value = 1.037502810292104
$ws.Range("J25").Value = 1.037389521592881
$ws.Range("K25").Value = 1.039505817631164
$ws.Range("L25").Value = 0.9944092447426414
$ws.Range("M25").Value = 1.032607065576322
$ws.Range("N25").Value = 1.016362481773597
